# qPCR_setup: update labware labels in the worktable to match the
# worklist worktable labware naming ("Tube" -> "Tubes[003]",
# "96-well"/"384-well" -> "96 Well[004]") and update the sheet
# selection accordingly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Column K on Sheet1 holds the "Sample labware" values for rows 5-127.
# Rows 5-40 were labelled "Tube" and are now "Tubes[003]".
# Rows 41-127 were labelled "96-well" / "384-well" and are now all
# unified under "96 Well[004]".
$ws1.Range("K5:K40").Value   = "Tubes[003]"
$ws1.Range("K41:K127").Value = "96 Well[004]"

# Update the active sheet's selection to match the edited area.
$ws1.Activate()
$ws1.Range("K41:K127").Select()
